$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D = 4 (Price), Column E = 5 (Volume(1h)), B = 2 (Coin), C = 3 (Link)
# Helper: write a value and force it to be stored as TEXT (not auto-coerced
# to a number), matching the source workbook's inlineStr cells, while
# resetting the cell style back to its original (unstyled) state afterwards
# so we don't leave stray NumberFormat/style residue behind.
function Set-TextValue($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue 2 4 "66.997.21"
$ws.Cells.Item(2, 5).Value = "  -1.84%  "

# Row 3 - Ethereum
Set-TextValue 3 4 "2.476.52"
$ws.Cells.Item(3, 5).Value = "  -2.44%  "

# Row 4 - TetherUSD
$ws.Cells.Item(4, 5).Value = "  -0.01%  "

# Row 5 - BNB
Set-TextValue 5 4 "585.61"
$ws.Cells.Item(5, 5).Value = "  -1.50%  "

# Row 6 - Solana
$ws.Cells.Item(6, 5).Value = "  -5.09%  "

# Row 7 - USDC
$ws.Cells.Item(7, 5).Value = "  +0.09%  "

# Row 8 - XRP
$ws.Cells.Item(8, 5).Value = "  -3.29%  "

# Row 9 - LidoStakedEther
Set-TextValue 9 4 "2.476.48"
$ws.Cells.Item(9, 5).Value = "  -2.40%  "

# Row 10 - Dogecoin
$ws.Cells.Item(10, 5).Value = "  -5.08%  "

# Row 11 - TRON
$ws.Cells.Item(11, 5).Value = "  +0.30%  "

# Row 12 - Toncoin
$ws.Cells.Item(12, 5).Value = "  -4.23%  "

# Row 13 - Cardano
Set-TextValue 13 4 "0.337"
$ws.Cells.Item(13, 5).Value = "  -3.37%  "

# Row 14 - Avalanche
Set-TextValue 14 4 "25.92"
$ws.Cells.Item(14, 5).Value = "  -4.12%  "

# Row 15 - WrappedliquidstakedEther2.0
Set-TextValue 15 4 "2.929.50"
$ws.Cells.Item(15, 5).Value = "  -2.30%  "

# Row 16 - ShibaInu
$ws.Cells.Item(16, 5).Value = "  -3.54%  "

# Row 17 - WrappedBTC
Set-TextValue 17 4 "66.855.79"

# Row 18 - WrappedEther
Set-TextValue 18 4 "2.468.21"
$ws.Cells.Item(18, 5).Value = "  -3.09%  "

# Row 19 - Chainlink
$ws.Cells.Item(19, 5).Value = "  +0.83%  "

# Row 20 - Uniswap
Set-TextValue 20 4 "7.88"
$ws.Cells.Item(20, 5).Value = "  -1.98%  "

# Row 21 - BitcoinCash
Set-TextValue 21 4 "362.56"
$ws.Cells.Item(21, 5).Value = "  -0.80%  "

# Row 22 - Polkadot
Set-TextValue 22 4 "4.10"
$ws.Cells.Item(22, 5).Value = "  -2.85%  "

# Row 23 - NEARProtocol
Set-TextValue 23 4 "4.45"
$ws.Cells.Item(23, 5).Value = "  -5.26%  "

# Row 24 - Dai
$ws.Cells.Item(24, 5).Value = "  +0.02%  "

# Row 25 - Litecoin
Set-TextValue 25 4 "70.84"
$ws.Cells.Item(25, 5).Value = "  -0.27%  "

# Row 26 - SuiNetwork
$ws.Cells.Item(26, 5).Value = "  -6.74%  "

# Row 27 - Aptos
$ws.Cells.Item(27, 5).Value = "  -8.11%  "

# Row 28 - Binance-PegBSC-USD
Set-TextValue 28 4 "0.998"
$ws.Cells.Item(28, 5).Value = "  +0.32%  "

# Row 29 - WrappedeETH
Set-TextValue 29 4 "2.615.55"
$ws.Cells.Item(29, 5).Value = "  -2.00%  "

# Row 30 - PEPE (subscript-3 character U+2083)
Set-TextValue 30 4 "0.0$([char]0x2083)0928"
$ws.Cells.Item(30, 5).Value = "  -7.06%  "

# Row 31 - InternetComputer(DFINITY)
$ws.Cells.Item(31, 5).Value = "  -2.27%  "

# Row 32 - Bittensor
Set-TextValue 32 4 "515.47"
$ws.Cells.Item(32, 5).Value = "  -6.50%  "

# Row 33 - PancakeSwap
$ws.Cells.Item(33, 5).Value = "  -2.72%  "

# Row 34 - Fetch.AI
$ws.Cells.Item(34, 5).Value = "  -6.21%  "

# Row 35 - FirstDigitalUSD
$ws.Cells.Item(35, 5).Value = "  -0.01%  "

# Row 36 - Kaspa
$ws.Cells.Item(36, 5).Value = "  -2.95%  "

# Row 37 - Monero
Set-TextValue 37 4 "158.12"
$ws.Cells.Item(37, 5).Value = "  +0.95%  "

# Row 38 - ImmutableX
Set-TextValue 38 4 "1.42"
$ws.Cells.Item(38, 5).Value = "  -3.99%  "

# Row 39 - EthereumClassic
Set-TextValue 39 4 "18.93"
$ws.Cells.Item(39, 5).Value = "  +0.13%  "

# Row 40 - WhiteBITCoin
Set-TextValue 40 4 "18.57"
$ws.Cells.Item(40, 5).Value = "  -0.61%  "

# Row 41 - Stacks
$ws.Cells.Item(41, 5).Value = "  -3.69%  "

# Row 42 - RenderToken
$ws.Cells.Item(42, 5).Value = "  -5.16%  "

# Row 43 - PolygonEcosystemToken
$ws.Cells.Item(43, 5).Value = "  -6.98%  "

# Row 44 - dogwifhat
$ws.Cells.Item(44, 5).Value = "  -2.99%  "

# Row 45 - OKB
Set-TextValue 45 4 "39.20"
$ws.Cells.Item(45, 5).Value = "  -2.23%  "

# Row 46 - Aave
Set-TextValue 46 4 "142.83"
$ws.Cells.Item(46, 5).Value = "  -3.28%  "

# Row 47 - ARBITRUM
$ws.Cells.Item(47, 5).Value = "  -4.96%  "

# Row 48 - Filecoin
$ws.Cells.Item(48, 5).Value = "  -3.88%  "

# Row 49 - BabyDogeCoin
$ws.Cells.Item(49, 5).Value = "  -4.50%  "

# Row 50 - Optimism
Set-TextValue 50 4 "1.65"

# Row 51 - Mantle -> Cronos
$ws.Cells.Item(51, 2).Value = "Cronos"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue 51 4 "0.0736"
$ws.Cells.Item(51, 5).Value = "  -2.84%  "
